$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (values like "1.00" or "63.011.24" must stay as literal text, not be
# reinterpreted as numbers/dates by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '62.924.55'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '3.382.35'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '563.35'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').Value = '154.51'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.386.06'
$ws.Range('E8').Value = '  +0.57%  '
$ws.Range('D9').Value = '0.541'
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('D12').Value = '0.431'
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('D13').Value = '3.975.33'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').Value = '0.0000187'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('D16').Value = '26.82'
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('D17').Value = '63.011.24'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '3.350.23'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '6.22'
$ws.Range('E19').Value = '  -4.68%  '
$ws.Range('D20').Value = '13.96'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').Value = '375.25'
$ws.Range('E21').Value = '  -3.72%  '
$ws.Range('D22').Value = '8.04'
$ws.Range('E22').Value = '  -5.22%  '
$ws.Range('D23').Value = '0.996'
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('D24').Value = '71.35'
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('D25').Value = '0.527'
$ws.Range('E25').Value = '  -3.10%  '
$ws.Range('E26').Value = '  +19.76%  '
$ws.Range('D27').Value = '9.39'
$ws.Range('E27').Value = '  +5.38%  '
$ws.Range('E28').Value = '  -2.89%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').Value = '5.98'
$ws.Range('E30').Value = '  +5.54%  '
$ws.Range('D31').Value = '1.98'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '22.95'
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '6.33'
$ws.Range('E35').Value = '  -5.20%  '
$ws.Range('D36').Value = '6.73'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').Value = '157.65'
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').Value = '1.44'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.0755'
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.905.20'
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('D41').Value = '1.81'
$ws.Range('E41').Value = '  -4.49%  '
$ws.Range('D42').Value = '26.62'
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').Value = '0.0314'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '41.06'
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.753'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '4.30'
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('D47').Value = '23.17'
$ws.Range('E47').Value = '  +4.13%  '
$ws.Range('D48').Value = '1.06'
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('D49').Value = '2.12'
$ws.Range('E49').Value = '  +15.32%  '
$ws.Range('D50').Value = '6.34'
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').Value = '0.830'
$ws.Range('E51').Value = '  +2.25%  '
